$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ncbitax")

# ---------------------------------------------------------------------------
# Column B width: widen to fit new longer vaccine-name strings
# (raw stored width closest achievable to target 19.140625)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.3

# ---------------------------------------------------------------------------
# Row 22: "any" / "Pandemrix" entry is reclassified as "pan" / "Pandemrix"
# and takes on the right-aligned "Arial" style (same as A11/A14, style 9)
# ---------------------------------------------------------------------------
$ws.Range("A11").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C24").Value = "Influenza A virus"
$ws.Range("A23").Value = "ut1"
$ws.Range("A24").Value = "ut2"

$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B23").Value = "unspecified-trivalent1"
$ws.Range("B24").Value = "unspecified-trivalent2"
$ws.Range("A22").Value = "pan"

# B22 keeps its existing "Pandemrix" text (unchanged)
$ws.Range("B22").Value = "Pandemrix"

# ---------------------------------------------------------------------------
# Row 23 (new row): unspecified trivalent vaccine #1
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = "H1N1 subtype"
$ws.Range("D23").Value = 114727
$ws.Range("E23").Value = "H3N2 subtype"
$ws.Range("F23").Value = 119210
$ws.Range("G23").Value = "Influenza B virus"
$ws.Range("H23").Value = 11520

# ---------------------------------------------------------------------------
# Row 24: replaced with unspecified trivalent vaccine #2 data
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 11320
$ws.Range("E24").Value = "Influenza B virus"
$ws.Range("F24").Value = 11520
$ws.Range("G24").ClearContents()
$ws.Range("H24").ClearContents()

# ---------------------------------------------------------------------------
# Update selection to reflect the extra rows now present on the sheet
# ---------------------------------------------------------------------------
$ws.Range("E29").Select()
